$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the previously-empty Department (E) and Institution (F) columns
# for the examiner rows with a "-" placeholder.
$ws.Range("E2").Value = "-"
$ws.Range("F2").Value = "-"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "-"
$ws.Range("E4").Value = "-"
$ws.Range("F4").Value = "-"
